$wb = $excel.ActiveWorkbook

# 1. Break the external link to data_designer.xlam (removes externalReferences + external link part)
$wb.BreakLink("data_designer.xlam", 1)

# 2. @core sheet: toggle show_locale_id_columns off (B4 True -> False)
$core = $wb.Worksheets.Item("@core")
$core.Range("B4").Value = $false

# 3. @core sheet: table "locale" header rename id:6 -> id:24
$core.Range("B9").Value = "id:24"

# 4. @examples sheet edits
$ex = $wb.Worksheets.Item("@examples")

# table_base data row: A3 0 -> "r"
$ex.Range("A3").Value = "r"

# table_ext: remove last column "Column1" (shrink range, clear leftover cells)
$loExt = $ex.ListObjects.Item("table_ext")
$loExt.Resize($ex.Range("A5:J6"))
$ex.Range("K5:K6").Clear()

# table_ext data row J6 JSON value change
$ex.Range("J6").Value = '{"Entry_0":["0","5557"]}'

# 5. Add new "modifiers" table (A8:I15) on @examples sheet
$ex.Range("A8").Value = "id:7"
$ex.Range("B8").Value = "label:label"
$ex.Range("C8").Value = "signature:support"
$ex.Range("D8").Value = "name:lid"
$ex.Range("E8").Value = "desc:lid"
$ex.Range("F8").Value = "note:lid"
$ex.Range("G8").Value = "name:ltext"
$ex.Range("H8").Value = "desc:ltext"
$ex.Range("I8").Value = "note:ltext"

$modData = @(
    @(0, "ENTITY_", 0, 1, 2),
    @(1, "ENTITY_1", 6, 7, 8),
    @(2, "ENTITY_2", 9, 10, 11),
    @(3, "ENTITY_3", 12, 13, 14),
    @(4, "ENTITY_4", 15, 16, 17),
    @(5, "ENTITY_5", 18, 19, 20),
    @(6, "ENTITY_6", 21, 22, 23)
)

$r = 9
foreach ($row in $modData) {
    $ex.Range("A$r").Value = $row[0]
    $ex.Range("B$r").Value = $row[1]
    $ex.Range("C$r").Formula = '=_xlfn.CONCAT(A' + $r + ',' + '" : "' + ',B' + $r + ')'
    $ex.Range("D$r").Value = $row[2]
    $ex.Range("E$r").Value = $row[3]
    $ex.Range("F$r").Value = $row[4]
    $ex.Range("G$r").Value = "Name"
    $ex.Range("H$r").Value = "Description"
    $ex.Range("I$r").Value = "Note"
    $r = $r + 1
}

Write-Host "stage2 ok"
